$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1137.4
$ws.Range("I29").Value = 91.5
$ws.Range("J29").Value = 1834.6666
$ws.Range("K29").Value = 274.5
$ws.Range("L29").Value = 5503.9998
$ws.Range("M29").Value = 6.5
$ws.Range("N29").Value = -6065.9998
$ws.Range("H38").Value = 501.33334
$ws.Range("I38").Value = 65
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 195
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 177
$ws.Range("N38").Value = -3744
$ws.Range("H43").Value = 89744760
$ws.Range("I43").Value = 200000320
$ws.Range("J43").Value = 20835040
$ws.Range("K43").Value = 200000320
$ws.Range("L43").Value = 20835040
$ws.Range("M43").Value = -200000251
$ws.Range("N43").Value = -20835178
$ws.Range("H58").Value = 491.5625
$ws.Range("I58").Value = 276.07144
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 828.21432
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -678.21432
$ws.Range("N58").Value = -6300
$ws.Range("H82").Value = 470.5
$ws.Range("I82").Value = 470.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1411.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1005.5
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 470.5
$ws.Range("I85").Value = 470.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1411.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -7.5
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 334062
$ws.Range("I86").Value = 500593
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 500593
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -499470
$ws.Range("N86").Value = -3246
$ws.Range("H89").Value = 334062
$ws.Range("I89").Value = 500593
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 2502965
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -2497349
$ws.Range("N89").Value = -16232

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2388.5
$ws.Range("I35").Value = 2191.3333
$ws.Range("J35").Value = 2980
$ws.Range("K35").Value = 2191.3333
$ws.Range("L35").Value = 2980
$ws.Range("M35").Value = -1785.3333
$ws.Range("N35").Value = -3792
$ws.Range("H36").Value = 50019.332
$ws.Range("I36").Value = 10000
$ws.Range("K36").Value = 10000
$ws.Range("M36").Value = -9654

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 7316.4287
$ws.Range("I36").Value = 1724.6666
$ws.Range("K36").Value = 1724.6666
$ws.Range("M36").Value = -1190.6666
$ws.Range("H92").Value = 31333.334
$ws.Range("J92").Value = 31333.334
$ws.Range("L92").Value = 31333.334
$ws.Range("N92").Value = -36325.334
$ws.Range("H95").Value = 28500
$ws.Range("J95").Value = 28500
$ws.Range("L95").Value = 28500
$ws.Range("N95").Value = -33992

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1978040.5
$ws.Range("I31").Value = 1990.2759
$ws.Range("J31").Value = 5348949.5
$ws.Range("K31").Value = 1990.2759
$ws.Range("L31").Value = 5348949.5
$ws.Range("M31").Value = -1695.2759
$ws.Range("N31").Value = -5349539.5
$ws.Range("H34").Value = 1978040.5
$ws.Range("I34").Value = 1990.2759
$ws.Range("J34").Value = 5348949.5
$ws.Range("K34").Value = 1990.2759
$ws.Range("L34").Value = 5348949.5
$ws.Range("M34").Value = -1788.2759
$ws.Range("N34").Value = -5349353.5
$ws.Range("H74").Value = 30733.334
$ws.Range("J74").Value = 30733.334
$ws.Range("L74").Value = 30733.334
$ws.Range("N74").Value = -32481.334
$ws.Range("H77").Value = 30733.334
$ws.Range("J77").Value = 30733.334
$ws.Range("L77").Value = 92200.00199999999
$ws.Range("N77").Value = -100936.002
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1029
$ws.Range("I134").Value = 774.6667
$ws.Range("J134").Value = 1283.3334
$ws.Range("K134").Value = 2324.0001
$ws.Range("L134").Value = 3850.0002
$ws.Range("M134").Value = 210.9998999999998
$ws.Range("N134").Value = -8920.0002
$ws.Range("H63").Value = 5959.923
$ws.Range("I63").Value = 3297.9
$ws.Range("J63").Value = 14833.333
$ws.Range("K63").Value = 9893.700000000001
$ws.Range("L63").Value = 44499.999
$ws.Range("M63").Value = -9144.700000000001
$ws.Range("N63").Value = -45997.999
$ws.Range("H66").Value = 5959.923
$ws.Range("I66").Value = 3297.9
$ws.Range("J66").Value = 14833.333
$ws.Range("K66").Value = 29681.1
$ws.Range("L66").Value = 133499.997
$ws.Range("M66").Value = -25937.1
$ws.Range("N66").Value = -140987.997
$ws.Range("H76").Value = 20147.143
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 20147.143
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 60441.429
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -61207.429
$ws.Range("H79").Value = 20147.143
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 20147.143
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 60441.429
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -63093.429
$ws.Range("H95").Value = 3000
$ws.Range("J95").Value = 3000
$ws.Range("L95").Value = 9000
$ws.Range("N95").Value = -13118
$ws.Range("H103").Value = 4073.8667
$ws.Range("I103").Value = 232.5
$ws.Range("J103").Value = 5470.727
$ws.Range("K103").Value = 697.5
$ws.Range("L103").Value = 16412.181
$ws.Range("M103").Value = 181.5
$ws.Range("N103").Value = -18170.181
$ws.Range("H114").Value = 1734
$ws.Range("I114").Value = 1184.5
$ws.Range("J114").Value = 2649.8333
$ws.Range("K114").Value = 3553.5
$ws.Range("L114").Value = 7949.499899999999
$ws.Range("M114").Value = -299.5
$ws.Range("N114").Value = -14457.4999
$ws.Range("H117").Value = 1601.1111
$ws.Range("I117").Value = 470
$ws.Range("J117").Value = 1742.5
$ws.Range("K117").Value = 1410
$ws.Range("L117").Value = 5227.5
$ws.Range("M117").Value = 2032
$ws.Range("N117").Value = -12111.5
$ws.Range("H121").Value = 1232186.4
$ws.Range("I121").Value = 322
$ws.Range("J121").Value = 1488824.8
$ws.Range("K121").Value = 966
$ws.Range("L121").Value = 4466474.4
$ws.Range("M121").Value = 344
$ws.Range("N121").Value = -4469094.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3759.5667
$ws.Range("I80").Value = 2863.5
$ws.Range("J80").Value = 6223.75
$ws.Range("K80").Value = 2863.5
$ws.Range("L80").Value = 6223.75
$ws.Range("M80").Value = -1865.5
$ws.Range("N80").Value = -8219.75
$ws.Range("H83").Value = 3759.5667
$ws.Range("I83").Value = 2863.5
$ws.Range("J83").Value = 6223.75
$ws.Range("K83").Value = 14317.5
$ws.Range("L83").Value = 31118.75
$ws.Range("M83").Value = -9325.5
$ws.Range("N83").Value = -41102.75
$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 716.875
$ws.Range("I22").Value = 733.3333
$ws.Range("J22").Value = 707
$ws.Range("K22").Value = 733.3333
$ws.Range("L22").Value = 707
$ws.Range("M22").Value = -438.3333
$ws.Range("N22").Value = -1297
$ws.Range("H27").Value = 716.875
$ws.Range("I27").Value = 733.3333
$ws.Range("J27").Value = 707
$ws.Range("K27").Value = 733.3333
$ws.Range("L27").Value = 707
$ws.Range("M27").Value = -626.3333
$ws.Range("N27").Value = -921
$ws.Range("H46").Value = 4000.1428
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3376

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 60195
$ws.Range("J57").Value = 60195
$ws.Range("L57").Value = 60195
$ws.Range("N57").Value = -61703
$ws.Range("H96").Value = 1166
$ws.Range("I96").Value = 749
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 749
$ws.Range("L96").Value = 2000
$ws.Range("N96").Value = -4746
$ws.Range("M96").Value = 624
